$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 0.4908867651578248
$ws.Cells.Item(2, 3).Value = 0.1396300584654213
$ws.Cells.Item(2, 4).Value = 0.07844862448544632
$ws.Cells.Item(2, 5).Value = 0.1054197549157543
$ws.Cells.Item(2, 7).Value = 0.6710566250718131
$ws.Cells.Item(2, 8).Value = 0.7869608080748733
$ws.Cells.Item(2, 9).Value = 0.7861401061115423
$ws.Cells.Item(2, 11).Value = 0.3203634580962103
$ws.Cells.Item(2, 12).Value = 0.2007034689628284
$ws.Cells.Item(2, 15).Value = 2.910641772078421
# Row 3
$ws.Cells.Item(3, 2).Value = 0.4466075507039022
$ws.Cells.Item(3, 3).Value = 0.1385154780506141
$ws.Cells.Item(3, 4).Value = 0.07117883871495678
$ws.Cells.Item(3, 5).Value = 0.1050550005728468
$ws.Cells.Item(3, 7).Value = 0.6746020438910563
$ws.Cells.Item(3, 8).Value = 0.7925544437119001
$ws.Cells.Item(3, 9).Value = 0.7932138554991539
$ws.Cells.Item(3, 11).Value = 0.2807909321530531
$ws.Cells.Item(3, 12).Value = 0.1931978078458201
$ws.Cells.Item(3, 15).Value = 2.929598318545388
# Row 4
$ws.Cells.Item(4, 2).Value = 0.4194760943727545
$ws.Cells.Item(4, 3).Value = 0.137834265489829
$ws.Cells.Item(4, 4).Value = 0.0667479500045971
$ws.Cells.Item(4, 5).Value = 0.1048858706590323
$ws.Cells.Item(4, 7).Value = 0.6772140653492258
$ws.Cells.Item(4, 8).Value = 0.7963240930669784
$ws.Cells.Item(4, 9).Value = 0.7979608947939951
$ws.Cells.Item(4, 11).Value = 0.2564374584504208
$ws.Cells.Item(4, 12).Value = 0.1886928588197918
$ws.Cells.Item(4, 15).Value = 2.942852278616201
# Row 5
$ws.Cells.Item(5, 2).Value = 0.4084346763355597
$ws.Cells.Item(5, 3).Value = 0.1375574884884614
$ws.Cells.Item(5, 4).Value = 0.06495060414337672
$ws.Cells.Item(5, 5).Value = 0.1048307566080133
$ws.Cells.Item(5, 7).Value = 0.6783878453149867
$ws.Cells.Item(5, 8).Value = 0.7979445879808793
$ws.Cells.Item(5, 9).Value = 0.7999968774837534
$ws.Cells.Item(5, 11).Value = 0.2464997996396647
$ws.Cells.Item(5, 12).Value = 0.1868831867703022
$ws.Cells.Item(5, 15).Value = 2.94865937373541
# Row 6
$ws.Cells.Item(6, 2).Value = 0.4066021775166746
$ws.Cells.Item(6, 3).Value = 0.1375115805090203
$ws.Cells.Item(6, 4).Value = 0.06465265710308188
$ws.Cells.Item(6, 5).Value = 0.1048224395081263
$ws.Cells.Item(6, 7).Value = 0.6785893540612236
$ws.Cells.Item(6, 8).Value = 0.7982187651481212
$ws.Cells.Item(6, 9).Value = 0.8003410829000934
$ws.Cells.Item(6, 11).Value = 0.2448488662914912
$ws.Cells.Item(6, 12).Value = 0.1865842730471741
$ws.Cells.Item(6, 15).Value = 2.949648158898057
# Row 7
$ws.Cells.Item(7, 2).Value = 0.419327124742324
$ws.Cells.Item(7, 3).Value = 0.137830529404539
$ws.Cells.Item(7, 4).Value = 0.06672367681254343
$ws.Cells.Item(7, 5).Value = 0.1048850714377068
$ws.Cells.Item(7, 7).Value = 0.6772294526346201
$ws.Cells.Item(7, 8).Value = 0.7963456060912932
$ws.Cells.Item(7, 9).Value = 0.7979879416891507
$ws.Cells.Item(7, 11).Value = 0.2563034892419012
$ws.Cells.Item(7, 12).Value = 0.1886683470241479
$ws.Cells.Item(7, 15).Value = 2.942928951255595
# Row 8
$ws.Cells.Item(8, 2).Value = 0.4756080705487591
$ws.Cells.Item(8, 3).Value = 0.1392451172850357
$ws.Cells.Item(8, 4).Value = 0.07593521409474135
$ws.Cells.Item(8, 5).Value = 0.1052826192063705
$ws.Cells.Item(8, 7).Value = 0.6721887427029571
$ws.Cells.Item(8, 8).Value = 0.7888199722924085
$ws.Cells.Item(8, 9).Value = 0.7884953439715972
$ws.Cells.Item(8, 11).Value = 0.3067307802190555
$ws.Cells.Item(8, 12).Value = 0.1980940755418743
$ws.Cells.Item(8, 15).Value = 2.916842886342707
# Row 9
$ws.Cells.Item(9, 2).Value = 0.586393108576118
$ws.Cells.Item(9, 3).Value = 0.1420428136741947
$ws.Cells.Item(9, 4).Value = 0.09425870889555199
$ws.Cells.Item(9, 5).Value = 0.1064965515868828
$ws.Cells.Item(9, 7).Value = 0.6657594119225791
$ws.Cells.Item(9, 8).Value = 0.7767188578870616
$ws.Cells.Item(9, 9).Value = 0.7730836907269314
$ws.Cells.Item(9, 11).Value = 0.4051540247912158
$ws.Cells.Item(9, 12).Value = 0.2173967324686714
$ws.Cells.Item(9, 15).Value = 2.878500438662797
# Row 10
$ws.Cells.Item(10, 2).Value = 0.6680133230510705
$ws.Cells.Item(10, 3).Value = 0.1441112595400043
$ws.Cells.Item(10, 4).Value = 0.1078800650700202
$ws.Cells.Item(10, 5).Value = 0.107652504325717
$ws.Cells.Item(10, 7).Value = 0.6631472907219802
$ws.Cells.Item(10, 8).Value = 0.7694445275403439
$ws.Cells.Item(10, 9).Value = 0.7637135928303138
$ws.Cells.Item(10, 11).Value = 0.4771599310628005
$ws.Cells.Item(10, 12).Value = 0.2320754425139171
$ws.Cells.Item(10, 15).Value = 2.858145149375417
# Row 11
$ws.Cells.Item(11, 2).Value = 0.7051881083759213
$ws.Cells.Item(11, 3).Value = 0.1450547731798579
$ws.Cells.Item(11, 4).Value = 0.114111552527163
$ws.Cells.Item(11, 5).Value = 0.1082355832309929
$ws.Cells.Item(11, 7).Value = 0.6624185739604229
$ws.Cells.Item(11, 8).Value = 0.7664855704099836
$ws.Cells.Item(11, 9).Value = 0.7598749622461476
$ws.Cells.Item(11, 11).Value = 0.5098465919952559
$ws.Cells.Item(11, 12).Value = 0.2388607429261924
$ws.Cells.Item(11, 15).Value = 2.850583009593777
# Row 12
$ws.Cells.Item(12, 2).Value = 0.7192711053929486
$ws.Cells.Item(12, 3).Value = 0.145412397541925
$ws.Cells.Item(12, 4).Value = 0.1164762818262801
$ws.Cells.Item(12, 5).Value = 0.108464593087028
$ws.Cells.Item(12, 7).Value = 0.6622087869078541
$ws.Cells.Item(12, 8).Value = 0.7654153926503682
$ws.Cells.Item(12, 9).Value = 0.7584823306062134
$ws.Cells.Item(12, 11).Value = 0.5222137118100534
$ws.Cells.Item(12, 12).Value = 0.2414456039127941
$ws.Cells.Item(12, 15).Value = 2.847963593781373
# Row 13
$ws.Cells.Item(13, 2).Value = 0.7162378387702972
$ws.Cells.Item(13, 3).Value = 0.1453353623049622
$ws.Cells.Item(13, 4).Value = 0.1159667733517438
$ws.Cells.Item(13, 5).Value = 0.1084149068571811
$ws.Cells.Item(13, 7).Value = 0.6622510244698105
$ws.Cells.Item(13, 8).Value = 0.7656436374248869
$ws.Cells.Item(13, 9).Value = 0.7587795465459735
$ws.Cells.Item(13, 11).Value = 0.5195507136244828
$ws.Cells.Item(13, 12).Value = 0.2408882238352987
$ws.Cells.Item(13, 15).Value = 2.848516869432785
# Row 14
$ws.Cells.Item(14, 2).Value = 0.7063466151440423
$ws.Cells.Item(14, 3).Value = 0.1450841886001228
$ws.Cells.Item(14, 4).Value = 0.1143060002612089
$ws.Cells.Item(14, 5).Value = 0.1082542595683087
$ws.Cells.Item(14, 7).Value = 0.6623999882138634
$ws.Cells.Item(14, 8).Value = 0.7663965180766183
$ws.Cells.Item(14, 9).Value = 0.7597591674595243
$ws.Cells.Item(14, 11).Value = 0.5108642587861709
$ws.Cells.Item(14, 12).Value = 0.2390730926210693
$ws.Cells.Item(14, 15).Value = 2.850362613978717
# Row 15
$ws.Cells.Item(15, 2).Value = 0.7002886750968571
$ws.Cells.Item(15, 3).Value = 0.1449303802680575
$ws.Cells.Item(15, 4).Value = 0.1132893787083731
$ws.Cells.Item(15, 5).Value = 0.1081569271493663
$ws.Cells.Item(15, 7).Value = 0.6624998514144664
$ws.Cells.Item(15, 8).Value = 0.7668642305625326
$ws.Cells.Item(15, 9).Value = 0.7603671548745368
$ws.Cells.Item(15, 11).Value = 0.5055421531775437
$ws.Cells.Item(15, 12).Value = 0.2379632769019366
$ws.Cells.Item(15, 15).Value = 2.851524990865528
# Row 16
$ws.Cells.Item(16, 2).Value = 0.6655847067108311
$ws.Cells.Item(16, 3).Value = 0.1440496475421753
$ws.Cells.Item(16, 4).Value = 0.1074735247257337
$ws.Cells.Item(16, 5).Value = 0.1076155487359109
$ws.Cells.Item(16, 7).Value = 0.6632041682135963
$ws.Cells.Item(16, 8).Value = 0.7696449453491283
$ws.Cells.Item(16, 9).Value = 0.7639729876710994
$ws.Cells.Item(16, 11).Value = 0.475022332190207
$ws.Cells.Item(16, 12).Value = 0.2316341696400883
$ws.Cells.Item(16, 15).Value = 2.85867352061959
# Row 17
$ws.Cells.Item(17, 2).Value = 0.6443059838626937
$ws.Cells.Item(17, 3).Value = 0.143509980552416
$ws.Cells.Item(17, 4).Value = 0.1039146378940359
$ws.Cells.Item(17, 5).Value = 0.1072980750901102
$ws.Cells.Item(17, 7).Value = 0.6637540024365052
$ws.Cells.Item(17, 8).Value = 0.771440483330494
$ws.Cells.Item(17, 9).Value = 0.7662936292377296
$ws.Cells.Item(17, 11).Value = 0.4562812196479058
$ws.Cells.Item(17, 12).Value = 0.2277790286446191
$ws.Cells.Item(17, 15).Value = 2.86349376997785
# Row 18
$ws.Cells.Item(18, 2).Value = 0.6320713410554788
$ws.Cells.Item(18, 3).Value = 0.1431998219668102
$ws.Cells.Item(18, 4).Value = 0.101870963291006
$ws.Cells.Item(18, 5).Value = 0.1071208605747351
$ws.Cells.Item(18, 7).Value = 0.6641135029039589
$ws.Cells.Item(18, 8).Value = 0.7725061914877642
$ws.Cells.Item(18, 9).Value = 0.7676682966135466
$ws.Cells.Item(18, 11).Value = 0.4454953545278215
$ws.Cells.Item(18, 12).Value = 0.2255718137422633
$ws.Cells.Item(18, 15).Value = 2.866426029117179
# Row 19
$ws.Cells.Item(19, 2).Value = 0.6279296670911663
$ws.Cells.Item(19, 3).Value = 0.1430948504445624
$ws.Cells.Item(19, 4).Value = 0.1011795787377139
$ws.Cells.Item(19, 5).Value = 0.107061784746648
$ws.Cells.Item(19, 7).Value = 0.6642426494418032
$ws.Cells.Item(19, 8).Value = 0.7728726844488563
$ws.Cells.Item(19, 9).Value = 0.7681405868267248
$ws.Cells.Item(19, 11).Value = 0.4418423543562255
$ws.Cells.Item(19, 12).Value = 0.2248262367167655
$ws.Cells.Item(19, 15).Value = 2.867446281081698
# Row 20
$ws.Cells.Item(20, 2).Value = 0.6465706994241884
$ws.Cells.Item(20, 3).Value = 0.1435674040523338
$ws.Cells.Item(20, 4).Value = 0.1042931459290628
$ws.Cells.Item(20, 5).Value = 0.107331313209464
$ws.Cells.Item(20, 7).Value = 0.6636909948152692
$ws.Cells.Item(20, 8).Value = 0.7712459340669682
$ws.Cells.Item(20, 9).Value = 0.7660424637635153
$ws.Cells.Item(20, 11).Value = 0.4582769177787895
$ws.Cells.Item(20, 12).Value = 0.2281883641491476
$ws.Cells.Item(20, 15).Value = 2.86296410898214
# Row 21
$ws.Cells.Item(21, 2).Value = 0.7092517586744123
$ws.Cells.Item(21, 3).Value = 0.1451579555419045
$ws.Cells.Item(21, 4).Value = 0.1147936742091673
$ws.Cells.Item(21, 5).Value = 0.1083012228926918
$ws.Cells.Item(21, 7).Value = 0.6623544377017652
$ws.Cells.Item(21, 8).Value = 0.7661740135540072
$ws.Cells.Item(21, 9).Value = 0.7594697740483767
$ws.Cells.Item(21, 11).Value = 0.5134159743197984
$ws.Cells.Item(21, 12).Value = 0.2396058230854692
$ws.Cells.Item(21, 15).Value = 2.849813845757495
# Row 22
$ws.Cells.Item(22, 2).Value = 0.7502503786686816
$ws.Cells.Item(22, 3).Value = 0.1461994175511165
$ws.Cells.Item(22, 4).Value = 0.121685501592907
$ws.Cells.Item(22, 5).Value = 0.108982963228847
$ws.Cells.Item(22, 7).Value = 0.6618666002180191
$ws.Cells.Item(22, 8).Value = 0.7631524770363427
$ws.Cells.Item(22, 9).Value = 0.7555295571380185
$ws.Cells.Item(22, 11).Value = 0.549390391223568
$ws.Cells.Item(22, 12).Value = 0.2471575968568231
$ws.Cells.Item(22, 15).Value = 2.842642828345475
# Row 23
$ws.Cells.Item(23, 2).Value = 0.7283658972040712
$ws.Cells.Item(23, 3).Value = 0.1456434025570559
$ws.Cells.Item(23, 4).Value = 0.1180045527475357
$ws.Cells.Item(23, 5).Value = 0.1086147338049024
$ws.Cells.Item(23, 7).Value = 0.662091652572073
$ws.Cells.Item(23, 8).Value = 0.7647383078582948
$ws.Cells.Item(23, 9).Value = 0.7575999973172358
$ws.Cells.Item(23, 11).Value = 0.5301960717050633
$ws.Cells.Item(23, 12).Value = 0.2431188885671389
$ws.Cells.Item(23, 15).Value = 2.846339860577473
# Row 24
$ws.Cells.Item(24, 2).Value = 0.6455468252639491
$ws.Cells.Item(24, 3).Value = 0.1435414425646186
$ws.Cells.Item(24, 4).Value = 0.1041220150401756
$ws.Cells.Item(24, 5).Value = 0.1073162697291039
$ws.Cells.Item(24, 7).Value = 0.6637193453991301
$ws.Cells.Item(24, 8).Value = 0.7713337856542921
$ws.Cells.Item(24, 9).Value = 0.7661558895219507
$ws.Cells.Item(24, 11).Value = 0.457374697958528
$ws.Cells.Item(24, 12).Value = 0.2280032750380485
$ws.Cells.Item(24, 15).Value = 2.863203067114682
# Row 25
$ws.Cells.Item(25, 2).Value = 0.5563810009588792
$ws.Cells.Item(25, 3).Value = 0.1412835786426641
$ws.Cells.Item(25, 4).Value = 0.08927383346765794
$ws.Cells.Item(25, 5).Value = 0.1061217243770001
$ws.Cells.Item(25, 7).Value = 0.6671282676534389
$ws.Cells.Item(25, 8).Value = 0.779708450567739
$ws.Cells.Item(25, 9).Value = 0.7769100127913724
$ws.Cells.Item(25, 11).Value = 0.3785799905148792
$ws.Cells.Item(25, 12).Value = 0.2120874097296905
$ws.Cells.Item(25, 15).Value = 2.887501006452396
